$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 596, shifting existing rows 596+ down by one
$ws.Rows.Item(596).Insert()

# Populate the new row 596 with the new record's data
$ws.Cells.Item(596, 1).Value = 10
$ws.Cells.Item(596, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(596, 3).Value = "La Araucanía"
$ws.Cells.Item(596, 4).Value = 45106
$ws.Cells.Item(596, 5).Value = 9
$ws.Cells.Item(596, 6).Value = 100112024
$ws.Cells.Item(596, 7).Value = "Choclo"
$ws.Cells.Item(596, 8).Value = "Dulce o Americano"
$ws.Cells.Item(596, 9).Value = "Primera"
$ws.Cells.Item(596, 10).Value = 145
$ws.Cells.Item(596, 11).Value = 20000
$ws.Cells.Item(596, 12).Value = 22000
$ws.Cells.Item(596, 13).Value = 20897
$ws.Cells.Item(596, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(596, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(596, 16).Value = 299
$ws.Cells.Item(596, 17).Value = 70
$ws.Cells.Item(596, 18).Value = "Hortaliza"
